$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their literal text formatting (Excel would
# otherwise coerce single-dot-separated values like "1.013" into numbers,
# or scientific notation for tiny values like "0.000008648").
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.923.67"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.842.54"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.011"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.57"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4764"
$ws.Range("E7").Value = "  +1.88%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3665"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07196"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9285"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.74"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07733"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.884.82"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.339"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.436"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.75"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.013"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008648"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.027.26"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.063"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.61"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.43"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.15"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.007"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.08"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.963"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08851"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.296"
$ws.Range("E31").Value = "  +4.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.173"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7395"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.488"
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.687"
$ws.Range("E35").Value = "  -5.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.110"
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01958"
$ws.Range("E37").Value = "  +1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05245"
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.957"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5230"
$ws.Range("E40").Value = "  +2.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.003"
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1510"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.256"
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.55"
$ws.Range("E44").Value = "  +2.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4727"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.012"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.58"
$ws.Range("E47").Value = "  +1.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.602"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.73"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06056"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8875"
$ws.Range("E51").Value = "  +3.34%  "
